$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, shifting all existing data rows
# (old rows 2..63) down to rows 3..64.
$ws.Rows("2:2").Insert()

# The inserted row inherited formatting from the row above (the header row);
# clear that so the new row matches the plain/default style used by every
# other data row, then restore the date format on column D.
$ws.Range("A2:T2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with the new data record.
$ws.Range("A2").Value = 8
$ws.Range("B2").Value = "Terminal La Palmera de La Serena"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44882
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103003
$ws.Range("J2").Value = "Damasco"
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 320
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 26000
$ws.Range("P2").Value = 25500
$ws.Range("Q2").Value = "$/caja 16 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1594
$ws.Range("T2").Value = 16
